# "Final project rough draft" — rename the two worksheets and switch the
# active tab/selection from Grade Cutoffs!B4 to Schedule!D22.

$wb = $excel.ActiveWorkbook

# Sheet1 -> Schedule, Sheet2 -> Grade Cutoffs
$wb.Worksheets.Item(1).Name = "Schedule"
$wb.Worksheets.Item(2).Name = "Grade Cutoffs"

# Make "Schedule" the active sheet and select cell D22 on it
$ws = $wb.Worksheets.Item("Schedule")
$ws.Activate()
$ws.Range("D22").Select()
